$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 = "OSMO_DEF", matching the style/formatting of the
# existing header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "OSMO_DEF"
